$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.404.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.684.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "683.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.682.17"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.90%  "

$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("E9").Value = "  -4.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.85%  "

$ws.Range("E12").Value = "  -7.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000232"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.309.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.688.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.406.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.39%  "

$ws.Range("E20").Value = "  -8.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.62%  "

$ws.Range("E23").Value = "  -8.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.828.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000123"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.22%  "

$ws.Range("E30").Value = "  -7.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.77%  "

$ws.Range("E32").Value = "  -8.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.658.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.78%  "

$ws.Range("E37").Value = "  -3.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -11.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.41%  "

$ws.Range("E42").Value = "  -8.87%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.941"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "165.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.35%  "

$ws.Range("E47").Value = "  -1.28%  "

$ws.Range("E48").Value = "  -12.12%  "

$ws.Range("E49").Value = "  -4.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000276"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "27.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.01%  "
